# Updates cryptos list data (prices & volume deltas) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.601.34"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "2.575.67"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'586.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'145.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").Value = "'0.106"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "3.035.85"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "63.479.52"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").Value = "2.579.80"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'11.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "'341.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'4.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'67.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").Value = "'1.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -3.61%  "
$ws.Range("D27").Value = "'7.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.40%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").Value = "'8.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'472.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0798"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").Value = "'176.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("D37").Value = "'18.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'158.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.87%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'40.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").Value = "'21.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").Value = "'18.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "'11.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  -3.72%  "
